$d = $word.ActiveDocument
$r = $d.Paragraphs(2).Range
$r.Find.Execute("Chapter#3", $true, $false, $false, $false, $false, $true, 1, $false, "Chapter#4", 2)
